$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting template to rows 1380-1406 (style block with F/G placeholders) ---
$srcFG = $ws.Range("A1370:H1370")
$dstFG = $ws.Range("A1380:H1406")
$srcFG.Copy()
$dstFG.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1380:A1406").EntireRow.RowHeight = 30

# --- Apply formatting template to rows 1407-1412 (style block without F/G) ---
foreach ($col in @("A","B","C","D","E","H")) {
  $s = $ws.Range($col + "1378")
  $d = $ws.Range($col + "1407:" + $col + "1412")
  $s.Copy()
  $d.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
$ws.Range("A1407:A1412").EntireRow.RowHeight = 30

# --- Write cell values ---
# Row 1379
$ws.Cells.Item(1379, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1379, 2).Value = '108538134'
$ws.Cells.Item(1379, 3).Value = '130TONG10BLK'
$ws.Cells.Item(1379, 4).Value = 'Visions 10 1/2" Black Disposable Plastic Tongs - 36/Case'
$ws.Cells.Item(1379, 5).Value = ' 2'
$ws.Cells.Item(1379, 8).Value = 34.49

# Row 1380
$ws.Cells.Item(1380, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1380, 2).Value = '108538134'
$ws.Cells.Item(1380, 3).Value = '409ML90360'
$ws.Cells.Item(1380, 4).Value = 'Mona Lisa Petit Four Marbled Chocolate Tulip Cup - 152/Box'
$ws.Cells.Item(1380, 5).Value = ' 1'
$ws.Cells.Item(1380, 8).Value = 186.49

# Row 1381
$ws.Cells.Item(1381, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1381, 2).Value = '108538134'
$ws.Cells.Item(1381, 3).Value = '176APRONVBLK'
$ws.Cells.Item(1381, 4).Value = 'Choice Black 38 Mil Heavy Weight Vinyl Dishwasher Apron - 40" x 25"'
$ws.Cells.Item(1381, 5).Value = ' 6'
$ws.Cells.Item(1381, 8).Value = 9.99

# Row 1382
$ws.Cells.Item(1382, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1382, 2).Value = '108538134'
$ws.Cells.Item(1382, 3).Value = '100CRAY4PKBX'
$ws.Cells.Item(1382, 4).Value = 'Choice 4 Pack Kids'' Restaurant Crayons in Print Box - 100/Case'
$ws.Cells.Item(1382, 5).Value = ' 1'
$ws.Cells.Item(1382, 8).Value = 20.99

# Row 1383
$ws.Cells.Item(1383, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1383, 2).Value = '108538134'
$ws.Cells.Item(1383, 3).Value = '40862028'
$ws.Cells.Item(1383, 4).Value = 'Ghirardelli 30 lb. Sweet Ground Chocolate & Cocoa Powder'
$ws.Cells.Item(1383, 5).Value = ' 3'
$ws.Cells.Item(1383, 8).Value = 123.47

# Row 1384
$ws.Cells.Item(1384, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1384, 2).Value = '108538134'
$ws.Cells.Item(1384, 3).Value = '104ADWHEYPRKT'
$ws.Cells.Item(1384, 4).Value = 'Add A Scoop Whey Protein Supplement Powder 2.5 lb. - 6/Case'
$ws.Cells.Item(1384, 5).Value = ' 1'
$ws.Cells.Item(1384, 8).Value = 316.99

# Row 1385
$ws.Cells.Item(1385, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1385, 2).Value = '108538134'
$ws.Cells.Item(1385, 3).Value = '40862029'
$ws.Cells.Item(1385, 4).Value = 'Ghirardelli 10 lb. Sweet Ground White Chocolate Flavored Powder'
$ws.Cells.Item(1385, 5).Value = ' 2'
$ws.Cells.Item(1385, 8).Value = 55.49

# Row 1386
$ws.Cells.Item(1386, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1386, 2).Value = '108538134'
$ws.Cells.Item(1386, 3).Value = '182RRF8'
$ws.Cells.Item(1386, 4).Value = 'Choice 6 1/2" x 7 3/4" Plastic Food Bag on a Roll - 2000/Case'
$ws.Cells.Item(1386, 5).Value = ' 6'
$ws.Cells.Item(1386, 8).Value = 19.29

# Row 1387
$ws.Cells.Item(1387, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1387, 2).Value = '108538134'
$ws.Cells.Item(1387, 3).Value = '711MAVYELLOW'
$ws.Cells.Item(1387, 4).Value = 'Mavalerio Yellow Sprinkles 25 lb.'
$ws.Cells.Item(1387, 5).Value = ' 1'
$ws.Cells.Item(1387, 8).Value = 59.99

# Row 1388
$ws.Cells.Item(1388, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1388, 2).Value = '108538134'
$ws.Cells.Item(1388, 3).Value = '3463318BK'
$ws.Cells.Item(1388, 4).Value = 'Fineline Platter Pleasers 3318-BK 8 1/4" Disposable Black Plastic Sandwich Spreader - 144/Case'
$ws.Cells.Item(1388, 5).Value = ' 2'
$ws.Cells.Item(1388, 8).Value = 38.49

# Row 1389
$ws.Cells.Item(1389, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1389, 2).Value = '108542341'
$ws.Cells.Item(1389, 3).Value = '245S12FU10R'
$ws.Cells.Item(1389, 4).Value = 'Enjay 1/2-10RS12 10" Fold-Under 1/2" Thick Silver Round Cake Drum - 12/Case'
$ws.Cells.Item(1389, 5).Value = ' 2'
$ws.Cells.Item(1389, 8).Value = 29.49

# Row 1390
$ws.Cells.Item(1390, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1390, 2).Value = '108542341'
$ws.Cells.Item(1390, 3).Value = '725CM5435'
$ws.Cells.Item(1390, 4).Value = 'Chefmaster 10.5 oz. Leaf Green Liqua-Gel Food Coloring'
$ws.Cells.Item(1390, 5).Value = ' 1'
$ws.Cells.Item(1390, 8).Value = 6.69

# Row 1391
$ws.Cells.Item(1391, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1391, 2).Value = '108542341'
$ws.Cells.Item(1391, 3).Value = '725CM5490'
$ws.Cells.Item(1391, 4).Value = 'Chefmaster 10.5 oz. Lemon Yellow Liqua-Gel Food Coloring'
$ws.Cells.Item(1391, 5).Value = ' 1'
$ws.Cells.Item(1391, 8).Value = 6.19

# Row 1392
$ws.Cells.Item(1392, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1392, 2).Value = '108542341'
$ws.Cells.Item(1392, 3).Value = '725CM5759'
$ws.Cells.Item(1392, 4).Value = 'Chefmaster 10.5 oz. Neon Brite Pink Liqua-Gel Food Coloring'
$ws.Cells.Item(1392, 5).Value = ' 6'
$ws.Cells.Item(1392, 8).Value = 11.99

# Row 1393
$ws.Cells.Item(1393, 1).Value = ' March 10, 2025'
$ws.Cells.Item(1393, 2).Value = '108542341'
$ws.Cells.Item(1393, 3).Value = '725CM5477'
$ws.Cells.Item(1393, 4).Value = 'Chefmaster 10.5 oz. Super Red Liqua-Gel Food Coloring'
$ws.Cells.Item(1393, 5).Value = ' 1'
$ws.Cells.Item(1393, 8).Value = 8.59

# Row 1394
$ws.Cells.Item(1394, 1).Value = ' March 11, 2025'
$ws.Cells.Item(1394, 2).Value = '108603925'
$ws.Cells.Item(1394, 3).Value = '580THERM'
$ws.Cells.Item(1394, 4).Value = 'Choice 4 1/2" Tube Refrigerator / Freezer Thermometer'
$ws.Cells.Item(1394, 5).Value = ' 10'
$ws.Cells.Item(1394, 8).Value = 1.69

# Row 1395
$ws.Cells.Item(1395, 1).Value = ' March 11, 2025'
$ws.Cells.Item(1395, 2).Value = '108603925'
$ws.Cells.Item(1395, 3).Value = '96515546'
$ws.Cells.Item(1395, 4).Value = 'Vitamix 15546 Drive Socket Set for BarBoss, Drink Machine, Blending Station, Portion Blending System, Vita-Prep, and Vita-Pro Series - 2/Pack'
$ws.Cells.Item(1395, 5).Value = ' 6'
$ws.Cells.Item(1395, 8).Value = 18.49

# Row 1396
$ws.Cells.Item(1396, 1).Value = ' March 11, 2025'
$ws.Cells.Item(1396, 2).Value = '108603925'
$ws.Cells.Item(1396, 3).Value = '9651195'
$ws.Cells.Item(1396, 4).Value = 'Vitamix 1195 64 oz. Clear Tritan™ Copolyester Blender Jar with Lid and Wet Blade Assembly'
$ws.Cells.Item(1396, 5).Value = ' 3'
$ws.Cells.Item(1396, 8).Value = 102.99

# Row 1397
$ws.Cells.Item(1397, 1).Value = ' March 11, 2025'
$ws.Cells.Item(1397, 2).Value = '108603925'
$ws.Cells.Item(1397, 3).Value = '96515978'
$ws.Cells.Item(1397, 4).Value = 'Vitamix 15978 Advance 48 oz. Clear Tritan™ Copolyester Deluxe Blender Jar with Lid and Wet Blade Assembly for Vitamix Blenders'
$ws.Cells.Item(1397, 5).Value = ' 2'
$ws.Cells.Item(1397, 8).Value = 116.99

# Row 1398
$ws.Cells.Item(1398, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1398, 2).Value = '108642953'
$ws.Cells.Item(1398, 3).Value = '588MILK632'
$ws.Cells.Item(1398, 4).Value = 'Urnex 12-MILK6-32 1 Liter Rinza Milk Frother Cleaner'
$ws.Cells.Item(1398, 5).Value = ' 12'
$ws.Cells.Item(1398, 8).Value = 17.99

# Row 1399
$ws.Cells.Item(1399, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1399, 2).Value = '108642953'
$ws.Cells.Item(1399, 3).Value = '544SYPFR013F'
$ws.Cells.Item(1399, 4).Value = 'Monin Premium Coconut Flavoring Syrup 1 Liter'
$ws.Cells.Item(1399, 5).Value = ' 8'
$ws.Cells.Item(1399, 8).Value = 10.49

# Row 1400
$ws.Cells.Item(1400, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1400, 2).Value = '108642953'
$ws.Cells.Item(1400, 3).Value = '544SYPFR034F'
$ws.Cells.Item(1400, 4).Value = 'Monin Premium Orange Flavoring / Fruit Syrup 1 Liter'
$ws.Cells.Item(1400, 5).Value = ' 8'
$ws.Cells.Item(1400, 8).Value = 10.49

# Row 1401
$ws.Cells.Item(1401, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1401, 2).Value = '108642953'
$ws.Cells.Item(1401, 3).Value = '544SYPFR042F'
$ws.Cells.Item(1401, 4).Value = 'Monin Premium Strawberry Flavoring / Fruit Syrup 1 Liter'
$ws.Cells.Item(1401, 5).Value = ' 8'
$ws.Cells.Item(1401, 8).Value = 10.49

# Row 1402
$ws.Cells.Item(1402, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1402, 2).Value = '108642953'
$ws.Cells.Item(1402, 3).Value = '544SYPFR063F'
$ws.Cells.Item(1402, 4).Value = 'Monin Premium White Chocolate Flavoring Syrup 1 Liter'
$ws.Cells.Item(1402, 5).Value = ' 8'
$ws.Cells.Item(1402, 8).Value = 10.49

# Row 1403
$ws.Cells.Item(1403, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1403, 2).Value = '108642953'
$ws.Cells.Item(1403, 3).Value = '544FRTRP066KT'
$ws.Cells.Item(1403, 4).Value = 'Monin 1 Liter Guava Fruit Puree - 4/Case'
$ws.Cells.Item(1403, 5).Value = ' 1'
$ws.Cells.Item(1403, 8).Value = 58.99

# Row 1404
$ws.Cells.Item(1404, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1404, 2).Value = '108647724'
$ws.Cells.Item(1404, 3).Value = '10200313'
$ws.Cells.Item(1404, 4).Value = 'Regal Chili Powder 5 lb.'
$ws.Cells.Item(1404, 5).Value = ' 1'
$ws.Cells.Item(1404, 8).Value = 20.99

# Row 1405
$ws.Cells.Item(1405, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1405, 2).Value = '108647724'
$ws.Cells.Item(1405, 3).Value = '245882WB'
$ws.Cells.Item(1405, 4).Value = '8" x 8" x 2 1/2" White Customizable Auto-Popup Window Bakery Box - 200/Bundle'
$ws.Cells.Item(1405, 5).Value = ' 1'
$ws.Cells.Item(1405, 8).Value = 72.49

# Row 1406
$ws.Cells.Item(1406, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1406, 2).Value = '108647724'
$ws.Cells.Item(1406, 3).Value = '245885CB'
$ws.Cells.Item(1406, 4).Value = '8" x 8" x 5" White Customizable Cake / Bakery Box - 100/Bundle'
$ws.Cells.Item(1406, 5).Value = ' 1'
$ws.Cells.Item(1406, 8).Value = 36.49

# Row 1407
$ws.Cells.Item(1407, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1407, 2).Value = '108650013'
$ws.Cells.Item(1407, 3).Value = '544SYPFR069F'
$ws.Cells.Item(1407, 4).Value = 'Monin Premium Blood Orange Flavoring Syrup 1 Liter'
$ws.Cells.Item(1407, 5).Value = ' 1'
$ws.Cells.Item(1407, 8).Value = 10.49

# Row 1408
$ws.Cells.Item(1408, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1408, 2).Value = '108650013'
$ws.Cells.Item(1408, 3).Value = '544SYPFR147F'
$ws.Cells.Item(1408, 4).Value = 'Monin Premium Elderflower Flavoring Syrup 1 Liter'
$ws.Cells.Item(1408, 5).Value = ' 1'
$ws.Cells.Item(1408, 8).Value = 10.49

# Row 1409
$ws.Cells.Item(1409, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1409, 2).Value = '108650013'
$ws.Cells.Item(1409, 3).Value = '544SYPFR095F'
$ws.Cells.Item(1409, 4).Value = 'Monin Premium Cucumber Flavoring Syrup 1 Liter'
$ws.Cells.Item(1409, 5).Value = ' 1'
$ws.Cells.Item(1409, 8).Value = 10.49

# Row 1410
$ws.Cells.Item(1410, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1410, 2).Value = '108650013'
$ws.Cells.Item(1410, 3).Value = '544SYPFR036F'
$ws.Cells.Item(1410, 4).Value = 'Monin Premium Peach Flavoring / Fruit Syrup 1 Liter'
$ws.Cells.Item(1410, 5).Value = ' 1'
$ws.Cells.Item(1410, 8).Value = 10.49

# Row 1411
$ws.Cells.Item(1411, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1411, 2).Value = '108650013'
$ws.Cells.Item(1411, 3).Value = '544SYPFR049F'
$ws.Cells.Item(1411, 4).Value = 'Monin Premium Granny Smith Apple Flavoring / Fruit Syrup 1 Liter'
$ws.Cells.Item(1411, 5).Value = ' 1'
$ws.Cells.Item(1411, 8).Value = 10.49

# Row 1412
$ws.Cells.Item(1412, 1).Value = ' March 12, 2025'
$ws.Cells.Item(1412, 2).Value = '108650013'
$ws.Cells.Item(1412, 3).Value = '544SYPAR047A'
$ws.Cells.Item(1412, 4).Value = 'Monin Premium Amaretto Flavoring Syrup 750 mL'
$ws.Cells.Item(1412, 5).Value = ' 1'
$ws.Cells.Item(1412, 8).Value = 8.29

# --- Row 1379 uses default/plain style (no custom formatting) ---
foreach ($col in @("A","B","C","D","E","H")) {
  $ws.Range($col + "1379").Style = "Normal"
}

# --- Update view selection to match ---
$ws.Range("A1380:G1412").Select()
